$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates for specific rows, per repull/recalculation of source data
$ws.Range("F2").Value  = -3
$ws.Range("F3").Value  = -1
$ws.Range("F7").Value  = -2
$ws.Range("F8").Value  = 3
$ws.Range("F11").Value = 4
$ws.Range("F18").Value = -1
$ws.Range("F20").Value = -1
$ws.Range("F23").Value = -2
$ws.Range("F27").Value = -1
$ws.Range("F30").Value = 1
$ws.Range("F40").Value = 1
$ws.Range("F41").Value = 6
